# Rewrites Sheet1 to the new "autofit columns" demo/readme layout:
#  - header row with a 4th "column name" label in col A
#  - instructions block in rows 2-3
#  - notes row (row 4)
#  - two sample data rows (rows 5-6)
#  - max/min length summary rows (rows 7-8), incl. a LEN/MAX formula in D7
#  - column widths roughly matched to the (manually autofit) content
#  - selection moved to column D (whole-column select)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear whatever was there before so stale cells don't linger ---
$ws.Cells.Clear()

# --- row 1: headers ---
$ws.Range("A1").Value = "column name"
$ws.Range("B1").Value = "Column 1"
$ws.Range("C1").Value = "Column 2"
$ws.Range("D1").Value = "Column 3"

# --- row 2: instructions banner (only col A) ---
$ws.Range("A2").Value = "INSTRUCTIONS BELOW"

# --- row 3: instructions detail row ---
$ws.Range("A3").Value = "in excel manually change column width to match  text in this row"
$ws.Range("B3").Value = "Colum"
$ws.Range("C3").Value = "short data"
$ws.Range("D3").Value = "C"

# --- row 4: notes row ---
$ws.Range("A4").Value = "notes"
$ws.Range("B4").Value = "longest col"
$ws.Range("C4").Value = "short data"
$ws.Range("D4").Value = "keep row 2 null"

# --- row 5: sample data row ---
$ws.Range("A5").Value = "data row"
$ws.Range("B5").Value = "longest row 2, column 1"
$ws.Range("C5").Value = "short data"

# --- row 6: sample data row ---
$ws.Range("A6").Value = "data row"
$ws.Range("B6").Value = "row 3"
$ws.Range("C6").Value = "row"
$ws.Range("D6").Value = "D"

# --- row 7: max length summary (D7 computed via formula) ---
$ws.Range("A7").Value = "max len in col"
$ws.Range("B7").Value = 28
$ws.Range("C7").Value = 22
$ws.Range("D7").Formula = "=MAX(LEN(D1), LEN(D5), LEN(D6), LEN(D4), LEN(D3))"

# --- row 8: min length summary ---
$ws.Range("A8").Value = "min len in col"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 0

# --- column widths (manually "autofit" to match the longest text per column) ---
$ws.Columns.Item(1).ColumnWidth = 18.333333333333332
$ws.Columns.Item(2).ColumnWidth = 5.333333333333333
$ws.Columns.Item(3).ColumnWidth = 9.166666666666666
$ws.Columns.Item(4).ColumnWidth = 0.8333333333333334

# --- selection parks on column D, matching the refreshed workbook state ---
$ws.Range("D:D").Select()

Write-Output "done"
